# Applies the crypto price/volume refresh described in the commit diff.
# For each changed cell: force text format (so numeric-looking strings like
# "1.00" / "0.0710" are preserved verbatim instead of being parsed as numbers),
# write the literal value, then restore the default "Normal" style so no
# stray per-cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "34.657.02" },
    @{ Cell = "E2"; Value = "  +1.55%  " },
    @{ Cell = "D3"; Value = "1.806.35" },
    @{ Cell = "E3"; Value = "  +0.81%  " },
    @{ Cell = "D4"; Value = "1.00" },
    @{ Cell = "E4"; Value = "  +0.00%  " },
    @{ Cell = "D5"; Value = "225.18" },
    @{ Cell = "E5"; Value = "  -1.02%  " },
    @{ Cell = "E6"; Value = "  -0.17%  " },
    @{ Cell = "E7"; Value = "  -0.03%  " },
    @{ Cell = "D8"; Value = "32.63" },
    @{ Cell = "E8"; Value = "  +3.98%  " },
    @{ Cell = "D9"; Value = "0.291" },
    @{ Cell = "E9"; Value = "  +3.33%  " },
    @{ Cell = "D10"; Value = "0.0710" },
    @{ Cell = "E10"; Value = "  +7.47%  " },
    @{ Cell = "D11"; Value = "0.0929" },
    @{ Cell = "E11"; Value = "  +0.23%  " },
    @{ Cell = "D12"; Value = "2.066.85" },
    @{ Cell = "E12"; Value = "  +0.86%  " },
    @{ Cell = "D13"; Value = "11.10" },
    @{ Cell = "E13"; Value = "  -3.65%  " },
    @{ Cell = "D14"; Value = "1.808.44" },
    @{ Cell = "E14"; Value = "  +0.95%  " },
    @{ Cell = "D15"; Value = "0.644" },
    @{ Cell = "E15"; Value = "  +0.85%  " },
    @{ Cell = "D16"; Value = "34.691.20" },
    @{ Cell = "E16"; Value = "  +1.64%  " },
    @{ Cell = "D17"; Value = "4.34" },
    @{ Cell = "E17"; Value = "  +2.35%  " },
    @{ Cell = "D18"; Value = "69.45" },
    @{ Cell = "E18"; Value = "  -0.24%  " },
    @{ Cell = "D19"; Value = "254.00" },
    @{ Cell = "E19"; Value = "  +0.14%  " },
    @{ Cell = "D20"; Value = "0.0₃0802" },
    @{ Cell = "E20"; Value = "  +7.69%  " },
    @{ Cell = "D21"; Value = "11.13" },
    @{ Cell = "E21"; Value = "  +6.07%  " },
    @{ Cell = "E22"; Value = "  -0.05%  " },
    @{ Cell = "E23"; Value = "  -0.71%  " },
    @{ Cell = "D24"; Value = "2.17" },
    @{ Cell = "E24"; Value = "  +1.44%  " },
    @{ Cell = "D25"; Value = "161.61" },
    @{ Cell = "E25"; Value = "  +2.81%  " },
    @{ Cell = "D26"; Value = "16.48" },
    @{ Cell = "E26"; Value = "  -0.93%  " },
    @{ Cell = "D27"; Value = "7.16" },
    @{ Cell = "E27"; Value = "  +1.44%  " },
    @{ Cell = "E28"; Value = "  +0.02%  " },
    @{ Cell = "D29"; Value = "684.21" },
    @{ Cell = "E29"; Value = "  +1,207.49%  " },
    @{ Cell = "D31"; Value = "0.0531" },
    @{ Cell = "E31"; Value = "  +2.72%  " },
    @{ Cell = "D32"; Value = "3.81" },
    @{ Cell = "E32"; Value = "  -0.39%  " },
    @{ Cell = "E33"; Value = "  -0.02%  " },
    @{ Cell = "D34"; Value = "3.64" },
    @{ Cell = "E34"; Value = "  +0.48%  " },
    @{ Cell = "D35"; Value = "1.89" },
    @{ Cell = "E35"; Value = "  +1.96%  " },
    @{ Cell = "D36"; Value = "1.439.71" },
    @{ Cell = "E36"; Value = "  -1.31%  " },
    @{ Cell = "B37"; Value = "ImmutableX" },
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" },
    @{ Cell = "D37"; Value = "0.647" },
    @{ Cell = "E37"; Value = "  +2.27%  " },
    @{ Cell = "B38"; Value = "TrustWalletToken" },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" },
    @{ Cell = "D38"; Value = "1.07" },
    @{ Cell = "E38"; Value = "  -0.48%  " },
    @{ Cell = "E39"; Value = "  +3.13%  " },
    @{ Cell = "D40"; Value = "85.17" },
    @{ Cell = "E40"; Value = "  +1.76%  " },
    @{ Cell = "D41"; Value = "0.960" },
    @{ Cell = "E41"; Value = "  +6.14%  " },
    @{ Cell = "E42"; Value = "  -1.02%  " },
    @{ Cell = "E43"; Value = "  +0.03%  " },
    @{ Cell = "D44"; Value = "2.17" },
    @{ Cell = "E44"; Value = "  +3.91%  " },
    @{ Cell = "D45"; Value = "6.07" },
    @{ Cell = "E45"; Value = "  +5.62%  " },
    @{ Cell = "E46"; Value = "  -1.00%  " },
    @{ Cell = "D47"; Value = "0.0496" },
    @{ Cell = "E47"; Value = "  -2.98%  " },
    @{ Cell = "D48"; Value = "1.959.93" },
    @{ Cell = "E48"; Value = "  +0.43%  " },
    @{ Cell = "E49"; Value = "  +2.22%  " },
    @{ Cell = "D50"; Value = "106.48" },
    @{ Cell = "E50"; Value = "  +8.67%  " },
    @{ Cell = "E51"; Value = "  +0.11%  " }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    $r.Style = "Normal"
}
